# Scheduled-runner style refresh of market/profit figures (columns H-N)
# across several Leve rows on each job sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 2997.5
$ws.Range("J76").Value = 2997
$ws.Range("L76").Value = 2997
$ws.Range("N76").Value = -3627

$ws.Range("H79").Value = 2997.5
$ws.Range("J79").Value = 2997
$ws.Range("L79").Value = 2997
$ws.Range("N79").Value = -5181

$ws.Range("H112").Value = 2540.2727
$ws.Range("J112").Value = 2694.3
$ws.Range("L112").Value = 8082.900000000001
$ws.Range("N112").Value = -10298.9

$ws.Range("H118").Value = 761.375
$ws.Range("I118").Value = 523.75
$ws.Range("J118").Value = 999
$ws.Range("K118").Value = 1571.25
$ws.Range("L118").Value = 2997
$ws.Range("M118").Value = 85.75
$ws.Range("N118").Value = -6311

$ws.Range("H132").Value = 1058.7333
$ws.Range("I132").Value = 1058.7333
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3176.199900000001
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -646.1999000000005
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 200
$ws.Range("I11").Value = 200
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 200
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -56
$ws.Range("N11").ClearContents()

$ws.Range("H63").Value = 2464.5454
$ws.Range("I63").Value = 710.4
$ws.Range("J63").Value = 20006
$ws.Range("K63").Value = 710.4
$ws.Range("L63").Value = 20006
$ws.Range("M63").Value = -24.39999999999998
$ws.Range("N63").Value = -21378

$ws.Range("H66").Value = 2464.5454
$ws.Range("I66").Value = 710.4
$ws.Range("J66").Value = 20006
$ws.Range("K66").Value = 3552
$ws.Range("L66").Value = 100030
$ws.Range("M66").Value = -120
$ws.Range("N66").Value = -106894

$ws.Range("H74").Value = 10691.637
$ws.Range("I74").Value = 12951.125
$ws.Range("J74").Value = 4666.3335
$ws.Range("K74").Value = 12951.125
$ws.Range("L74").Value = 4666.3335
$ws.Range("M74").Value = -12077.125
$ws.Range("N74").Value = -6414.3335

$ws.Range("H77").Value = 10691.637
$ws.Range("I77").Value = 12951.125
$ws.Range("J77").Value = 4666.3335
$ws.Range("K77").Value = 64755.625
$ws.Range("L77").Value = 23331.6675
$ws.Range("M77").Value = -60387.625
$ws.Range("N77").Value = -32067.6675

$ws.Range("H107").Value = 170000
$ws.Range("J107").Value = 170000
$ws.Range("L107").Value = 170000
$ws.Range("N107").Value = -177680

$ws.Range("H132").Value = 2493.4285
$ws.Range("I132").Value = 1961.8572
$ws.Range("J132").Value = 3556.5715
$ws.Range("K132").Value = 5885.571599999999
$ws.Range("L132").Value = 10669.7145
$ws.Range("M132").Value = -3355.571599999999
$ws.Range("N132").Value = -15729.7145

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2010.3334
$ws.Range("I134").Value = 2010.3334
$ws.Range("K134").Value = 6031.0002
$ws.Range("M134").Value = -3496.0002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 1502.5
$ws.Range("J13").Value = 1502.5
$ws.Range("L13").Value = 1502.5
$ws.Range("N13").Value = -1780.5

$ws.Range("H16").Value = 1699.6
$ws.Range("I16").Value = 1749.5
$ws.Range("K16").Value = 1749.5
$ws.Range("M16").Value = -1462.5

$ws.Range("H31").Value = 2341.9412
$ws.Range("J31").Value = 2683
$ws.Range("L31").Value = 2683
$ws.Range("N31").Value = -3273

$ws.Range("H34").Value = 2341.9412
$ws.Range("J34").Value = 2683
$ws.Range("L34").Value = 2683
$ws.Range("N34").Value = -3087

$ws.Range("H105").Value = 4783.7144
$ws.Range("I105").Value = 5123.25
$ws.Range("K105").Value = 5123.25
$ws.Range("M105").Value = -3376.25

$ws.Range("H107").Value = 1255.5714
$ws.Range("I107").Value = 1457.9333
$ws.Range("K107").Value = 1457.9333
$ws.Range("M107").Value = 462.0667000000001

$ws.Range("H113").Value = 1699.6
$ws.Range("I113").Value = 1749.5
$ws.Range("K113").Value = 1749.5
$ws.Range("M113").Value = 420.5

$ws.Range("H132").Value = 2228.8125
$ws.Range("I132").Value = 1796.4546
$ws.Range("J132").Value = 3180
$ws.Range("K132").Value = 5389.3638
$ws.Range("L132").Value = 9540
$ws.Range("M132").Value = -2859.3638
$ws.Range("N132").Value = -14600

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 416
$ws.Range("I12").Value = 600
$ws.Range("K12").Value = 1800
$ws.Range("M12").Value = -1627

$ws.Range("H56").Value = 9866.4
$ws.Range("I56").Value = 9866.4
$ws.Range("K56").Value = 9866.4
$ws.Range("M56").Value = -9336.4

$ws.Range("H98").Value = 1681.4
$ws.Range("J98").Value = 738.25
$ws.Range("L98").Value = 2214.75
$ws.Range("N98").Value = -5210.75

$ws.Range("H113").Value = 1261.5
$ws.Range("I113").Value = 692
$ws.Range("K113").Value = 2076
$ws.Range("M113").Value = 94

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 6511.6665
$ws.Range("J102").Value = 8885.5
$ws.Range("L102").Value = 8885.5
$ws.Range("N102").Value = -12129.5

$ws.Range("H132").Value = 4031.5715
$ws.Range("I132").Value = 3644.6
$ws.Range("K132").Value = 10933.8
$ws.Range("M132").Value = -8403.799999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2628.7058
$ws.Range("I7").Value = 2646.0667
$ws.Range("K7").Value = 2646.0667
$ws.Range("M7").Value = -2534.0667

$ws.Range("H122").Value = 4362.8335
$ws.Range("I122").Value = 4460.9414
$ws.Range("J122").Value = 2695
$ws.Range("K122").Value = 13382.8242
$ws.Range("L122").Value = 8085
$ws.Range("M122").Value = -10932.8242
$ws.Range("N122").Value = -12985

$ws.Range("H126").Value = 2628.7058
$ws.Range("I126").Value = 2646.0667
$ws.Range("K126").Value = 7938.2001
$ws.Range("M126").Value = -5468.2001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5078.5
$ws.Range("J62").Value = 5594.2
$ws.Range("L62").Value = 5594.2
$ws.Range("N62").Value = -6842.2

$ws.Range("H65").Value = 5078.5
$ws.Range("J65").Value = 5594.2
$ws.Range("L65").Value = 27971
$ws.Range("N65").Value = -34211

$ws.Range("H126").Value = 1671.5
$ws.Range("I126").Value = 1650.125
$ws.Range("K126").Value = 4950.375
$ws.Range("M126").Value = -2480.375

$ws.Range("H136").Value = 5217.2856
$ws.Range("I136").Value = 4334.5
$ws.Range("J136").Value = 7424.25
$ws.Range("K136").Value = 13003.5
$ws.Range("L136").Value = 22272.75
$ws.Range("M136").Value = -10453.5
$ws.Range("N136").Value = -27372.75
